$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', RobustScaler()), ('selector', 'passthrough'),`n                ('model',`n                 BaggingClassifier(estimator=RandomForestClassifier(max_depth=6,`n                                                                    min_samples_leaf=6,`n                                                                    min_samples_split=3,`n                                                                    n_estimators=5,`n                                                                    random_state=42),`n                                   random_state=42))])"
$ws.Range("B2").Value = 0.6761904761904762
$ws.Range("C2").Value = "{'scaler': RobustScaler(), 'model__n_estimators': 10, 'model__estimator__n_estimators': 5, 'model__estimator__min_samples_split': 3, 'model__estimator__min_samples_leaf': 6, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 6, 'model__estimator__class_weight': None}"
$ws.Range("D2").Value = 0.5000000000000001
$ws.Range("E2").Value = "[1 0 0 1 0 0 1 1 0 1 0 0]"
$ws.Range("F2").Value = "[0 1 1 1 1 1 1 1 1 1 1 1]"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.8130238095238096
$ws.Range("I2").Value = 0.02263941009981263
$ws.Range("J2").Value = 0.5887619047619046
$ws.Range("K2").Value = 0.06650037095434869

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', None), ('selector', 'passthrough'),`n                ('model',`n                 BaggingClassifier(estimator=RandomForestClassifier(class_weight='balanced',`n                                                                    max_depth=2,`n                                                                    min_samples_leaf=6,`n                                                                    n_estimators=50,`n                                                                    random_state=42),`n                                   random_state=42))])"
$ws.Range("B3").Value = 0.6
$ws.Range("C3").Value = "{'scaler': None, 'model__n_estimators': 10, 'model__estimator__n_estimators': 50, 'model__estimator__min_samples_split': 2, 'model__estimator__min_samples_leaf': 6, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 2, 'model__estimator__class_weight': 'balanced'}"
$ws.Range("D3").Value = 0.5
$ws.Range("E3").Value = "[1 0 1 0 0 0 0 1 1 0 1 1]"
$ws.Range("F3").Value = "[1 1 0 0 1 0 1 1 0 0 1 0]"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.8085714285714286
$ws.Range("I3").Value = 0.02757051045080954
$ws.Range("J3").Value = 0.5472380952380952
$ws.Range("K3").Value = 0.05118580728890872

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None), ('selector', 'passthrough'),`n                ('model',`n                 BaggingClassifier(estimator=RandomForestClassifier(class_weight='balanced',`n                                                                    max_depth=2,`n                                                                    min_samples_leaf=4,`n                                                                    min_samples_split=5,`n                                                                    n_estimators=10,`n                                                                    random_state=42),`n                                   random_state=42))])"
$ws.Range("B4").Value = 0.6
$ws.Range("C4").Value = "{'scaler': None, 'model__n_estimators': 10, 'model__estimator__n_estimators': 10, 'model__estimator__min_samples_split': 5, 'model__estimator__min_samples_leaf': 4, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 2, 'model__estimator__class_weight': 'balanced'}"
$ws.Range("D4").Value = 0.7777777777777777
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[1 1 1 1 1 1 0 0 1 1 1 1]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.8514047619047619
$ws.Range("I4").Value = 0.02418291686109236
$ws.Range("J4").Value = 0.5278095238095237
$ws.Range("K4").Value = 0.07201148000967575
